# Update "想去人数" (interest count) figures in the 展览 and 全部类型 sheets
# to reflect freshly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 660
    "F6"  = 47
    "F7"  = 41
    "F8"  = 2088
    "F9"  = 4133
    "F10" = 98
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
